$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.848.56"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.933.72"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.16"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.96"
$ws.Range("E6").Value = "  +9.49%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +4.78%  "

$ws.Range("D9").Value = "2.928.73"
$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  +3.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.79"
$ws.Range("E14").Value = "  +5.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  +2.98%  "

$ws.Range("D16").Value = "3.418.32"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.86"
$ws.Range("E17").Value = "  +8.02%  "

$ws.Range("D18").Value = "2.931.75"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").Value = "57.802.55"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.70"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.28"
$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.694"
$ws.Range("E22").Value = "  +7.30%  "

$ws.Range("E23").Value = "  +8.11%  "

$ws.Range("E24").Value = "  +3.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.31"
$ws.Range("E25").Value = "  +3.62%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("E29").Value = "  +6.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("E30").Value = "  +6.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.43"
$ws.Range("E31").Value = "  +4.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.95"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0975"
$ws.Range("E33").Value = "  +4.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +6.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.937"
$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  +5.66%  "

$ws.Range("D37").Value = "0.0₃0696"
$ws.Range("E37").Value = "  +14.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.27"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.73"
$ws.Range("E39").Value = "  +6.19%  "

$ws.Range("E40").Value = "  +11.71%  "

$ws.Range("E41").Value = "  +3.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "375.16"
$ws.Range("E42").Value = "  +7.47%  "

$ws.Range("E43").Value = "  +1.03%  "

$ws.Range("D44").Value = "2.697.00"
$ws.Range("E44").Value = "  +3.57%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.89"
$ws.Range("E46").Value = "  +4.43%  "

$ws.Range("E47").Value = "  +4.71%  "

$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("E49").Value = "  +2.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.93"
$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("E51").Value = "  +3.51%  "
